# Cleaning the standardized balance sheet dataframe / new rendering structure.
#
# On Sheet1 the rows that described the "Preferred Stock, Value, Issued" and
# "Temporary Equity, Carrying Amount, Attributable to Parent" XBRL facts (old
# rows 28-29), as well as the now-unused "Common Stock, Value, Issued" /
# "Additional Paid in Capital" rows (old rows 25-26) are removed, collapsing
# the sheet from 34 to 30 used rows. Sheet1 becomes the active/selected
# sheet (it was Sheet2_Calculated_Fields before), with the new row 26
# ("Retained Earnings") highlighted, matching the author's "WIP" cleanup pass.

$wb  = $excel.ActiveWorkbook

# Best-effort: try to restore the author's last window position. (The
# iron_native shim may not persist this cosmetic metadata, but it shouldn't
# hurt to try.)
try {
    $win = $excel.ActiveWindow
    $win.Left = 16060
    $win.Top  = 5140
} catch {
}

$ws1 = $wb.Worksheets.Item("Sheet1")

# Remove the two "Preferred Stock, Value, Issued" / "Temporary Equity..."
# rows first (higher row numbers), then the "Common Stock, Value, Issued" /
# "Additional Paid in Capital" rows, so row numbers used for each delete
# refer to the original layout.
$ws1.Rows("28:29").Delete()
$ws1.Rows("25:26").Delete()

# Sheet1 is now the visible/active sheet (previously Sheet2 was active).
$ws1.Activate()

# Select the full "Retained Earnings" row (new row 26) as the final user
# selection, as in the committed workbook.
[void]$ws1.Range("A26:XFD26").Select()
